$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename row 14 Name from "GTA: San Andreas" to "Grand Theft Auto: San Andreas (GTA)"
# (all other fields in row 14 stay the same)

# Row 16: Resident Evil 2
$ws.Cells.Item(16, 1).Value = 'Resident Evil 2'
$ws.Cells.Item(16, 2).Value = 'January 25, 2019'
$ws.Cells.Item(16, 3).Value = 91
$ws.Cells.Item(16, 4).Value = 'Action, Shooter, Adventure'
$ws.Cells.Item(16, 5).Value = 'PC, Playstation, Xbox'

# Row 17: Ori and the Will of the Wisps
$ws.Cells.Item(17, 1).Value = 'Ori and the Will of the Wisps'
$ws.Cells.Item(17, 2).Value = 'March 10, 2020'
$ws.Cells.Item(17, 3).Value = 91
$ws.Cells.Item(17, 4).Value = 'Action, Adventure, Platformer'
$ws.Cells.Item(17, 5).Value = 'PC Xbox, Nintendo'

# Row 18: The Last of Us Part II
$ws.Cells.Item(18, 1).Value = 'The Last of Us Part II'
$ws.Cells.Item(18, 2).Value = 'June 19, 2020'
$ws.Cells.Item(18, 3).Value = 93
$ws.Cells.Item(18, 4).Value = 'Action, Shooter, Adventure'
$ws.Cells.Item(18, 5).Value = 'Playstation'

# Row 19: Max Payne
$ws.Cells.Item(19, 1).Value = 'Max Payne'
$ws.Cells.Item(19, 2).Value = 'July 23, 2001'
$ws.Cells.Item(19, 3).Value = 89
$ws.Cells.Item(19, 4).Value = 'Action, Shooter'
$ws.Cells.Item(19, 5).Value = 'PC, Playstation, Xbox'

# Row 20: Ghost of Tsushima
$ws.Cells.Item(20, 1).Value = 'Ghost of Tsushima'
$ws.Cells.Item(20, 2).Value = 'July 17, 2020'
$ws.Cells.Item(20, 3).Value = 83
$ws.Cells.Item(20, 4).Value = 'Action, Adventure, RPG'
$ws.Cells.Item(20, 5).Value = 'Playstation'

# Row 21: Marvel's Spider-Man
$ws.Cells.Item(21, 1).Value = 'Marvel''s Spider-Man'
$ws.Cells.Item(21, 2).Value = 'September 17, 2020'
$ws.Cells.Item(21, 3).Value = 87
$ws.Cells.Item(21, 4).Value = 'Action, Adventure'
$ws.Cells.Item(21, 5).Value = 'PC, Playstation'

# Row 22: Metal Gear Solid
$ws.Cells.Item(22, 1).Value = 'Metal Gear Solid'
$ws.Cells.Item(22, 2).Value = 'September 13, 1998'
$ws.Cells.Item(22, 3).Value = 94
$ws.Cells.Item(22, 4).Value = 'Action, Shooter, Adventure'
$ws.Cells.Item(22, 5).Value = 'PC, Playstation'

# Row 23: Warcraft 3: Reign of Chaos
$ws.Cells.Item(23, 1).Value = 'Warcraft 3: Reign of Chaos'
$ws.Cells.Item(23, 2).Value = 'June 1, 2002'
$ws.Cells.Item(23, 3).Value = 92
$ws.Cells.Item(23, 4).Value = 'Strategy'
$ws.Cells.Item(23, 5).Value = 'PC'

# Row 24: Silent Hill 2
$ws.Cells.Item(24, 1).Value = 'Silent Hill 2'
$ws.Cells.Item(24, 2).Value = 'September 24, 2001'
$ws.Cells.Item(24, 3).Value = 89
$ws.Cells.Item(24, 4).Value = 'Action, Adventure'
$ws.Cells.Item(24, 5).Value = 'PC, Playstation, Xbox'

# Row 25: Super Mario Odyssey
$ws.Cells.Item(25, 1).Value = 'Super Mario Odyssey'
$ws.Cells.Item(25, 2).Value = 'October 27, 2017'
$ws.Cells.Item(25, 3).Value = 97
$ws.Cells.Item(25, 4).Value = 'Arcade, Platformer'
$ws.Cells.Item(25, 5).Value = 'Nintendo'

# Row 14 rename (Name only)
$ws.Cells.Item(14, 1).Value = 'Grand Theft Auto: San Andreas (GTA)'

# Row 26: Grand Theft Auto V (GTA)
$ws.Cells.Item(26, 1).Value = 'Grand Theft Auto V (GTA)'
$ws.Cells.Item(26, 2).Value = 'September 17. 2013'
$ws.Cells.Item(26, 3).Value = 92
$ws.Cells.Item(26, 4).Value = 'Action, Adventure'
$ws.Cells.Item(26, 5).Value = 'PC, Playstation, Xbox'

# Row 27: Half-Life: Alyx
$ws.Cells.Item(27, 1).Value = 'Half-Life: Alyx'
$ws.Cells.Item(27, 2).Value = 'March 23, 2020'
$ws.Cells.Item(27, 3).Value = 93
$ws.Cells.Item(27, 4).Value = 'Action, Shooter, Adventure'
$ws.Cells.Item(27, 5).Value = 'PC'

# Row 28: Elden Ring
$ws.Cells.Item(28, 1).Value = 'Elden Ring'
$ws.Cells.Item(28, 2).Value = 'February 25, 2022'
$ws.Cells.Item(28, 3).Value = 95
$ws.Cells.Item(28, 4).Value = 'Action, RPG'
$ws.Cells.Item(28, 5).Value = 'PC, Playstation, Xbox'

# Row 29: Final Fantasy VII (1997)
$ws.Cells.Item(29, 1).Value = 'Final Fantasy VII (1997)'
$ws.Cells.Item(29, 2).Value = 'January 31, 1997'
$ws.Cells.Item(29, 4).Value = 'Action, Adventure, RPG'
$ws.Cells.Item(29, 5).Value = 'PC, Playstation, Xbox'

# Row 30: Diablo II
$ws.Cells.Item(30, 1).Value = 'Diablo II'
$ws.Cells.Item(30, 2).Value = 'June 29, 2000'
$ws.Cells.Item(30, 3).Value = 88
$ws.Cells.Item(30, 4).Value = 'Action, RPG'
$ws.Cells.Item(30, 5).Value = 'PC'

# Row 31: Need For Speed: Most Wanted
$ws.Cells.Item(31, 1).Value = 'Need For Speed: Most Wanted'
$ws.Cells.Item(31, 2).Value = 'Nov 15, 2005'
$ws.Cells.Item(31, 3).Value = 83
$ws.Cells.Item(31, 4).Value = 'Racing, Arcade'
$ws.Cells.Item(31, 5).Value = 'PC, Playstation, Xbox, Nintendo'

# Row 32: Metal Gear Solid 3: Snake Eater
$ws.Cells.Item(32, 1).Value = 'Metal Gear Solid 3: Snake Eater'
$ws.Cells.Item(32, 2).Value = 'November 17, 2004'
$ws.Cells.Item(32, 3).Value = 84
$ws.Cells.Item(32, 4).Value = 'Action'
$ws.Cells.Item(32, 5).Value = 'Playstation, Xbox, Nintendo'

# Row 33: Heroes of Might and Magic 3: The Restoration of Erathia
$ws.Cells.Item(33, 1).Value = 'Heroes of Might and Magic 3: The Restoration of Erathia'
$ws.Cells.Item(33, 2).Value = 'March 3, 1999'
$ws.Cells.Item(33, 4).Value = 'Strategy'
$ws.Cells.Item(33, 5).Value = 'PC'

# Row 34: Fallout: New Vegas
$ws.Cells.Item(34, 1).Value = 'Fallout: New Vegas'
$ws.Cells.Item(34, 2).Value = 'October 19, 2010'
$ws.Cells.Item(34, 3).Value = 84
$ws.Cells.Item(34, 4).Value = 'Action, Shooter, Adventure, RPG'
$ws.Cells.Item(34, 5).Value = 'PC, Playstation, Xbox'

# Row 35/36 were filled slightly out of order by the original author (date
# columns and row 36 before row 35's name) - replicate that exact sequence
# so new entries land in the same slots.
$ws.Cells.Item(35, 2).Value = 'October 27, 2002'
$ws.Cells.Item(35, 3).Value = 94
$ws.Cells.Item(35, 4).Value = 'Action, Adventure'
$ws.Cells.Item(35, 5).Value = 'PC, Playstation, Xbox'

# Row 36: Shadow of the Colossus
$ws.Cells.Item(36, 1).Value = 'Shadow of the Colossus'
$ws.Cells.Item(36, 2).Value = 'October 18, 2005'
$ws.Cells.Item(36, 3).Value = 91
$ws.Cells.Item(36, 4).Value = 'Action, Adventure, RPG'
$ws.Cells.Item(36, 5).Value = 'Playstation'

# Row 35 name filled in last
$ws.Cells.Item(35, 1).Value = 'Grand Theft Auto: Vice City (GTA)'

# Autofit columns to reflect the new, longer content (column A holds game
# names which are now noticeably longer, column D holds platform lists)
$ws.Columns("A:A").AutoFit()
$ws.Columns("D:D").AutoFit()

# Move the active selection to the first empty row below the data, same as
# Excel leaves it after the last row of data is typed in
$ws.Range("A37").Select()
